# Add a new worksheet ("Sheet2") right after the existing "Sheet1",
# put some content in it, and leave it as the active sheet/selection
# (mirrors a user adding a sheet, typing in A1, and pressing Enter).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Range("A1").Value = "This is sheet 2"
$ws2.Range("A2").Select() | Out-Null
